# Generate Report for Handback
# Updates the "6daa4734-0ac8-48c1-b198-78425699ee4e" row (row 6) on the
# zh-cn and de-de sheets: it now has an out-of-date handback, so we stamp a
# handback datetime, link/point at the produced target xliff, and record the
# "not latest version" error detail - plus widen the Error Detail column so
# the message is readable.

$wb = $excel.ActiveWorkbook

$ColWidth = 40 - (5 / 6)   # COM ColumnWidth units -> 40 chars in the saved XML

$rows = @(
    @{
        SheetName   = "zh-cn"
        TargetFile  = "6daa4734-0ac8-48c1-b198-78425699ee4e.8db02af582f18e84bbda947aeb8bdd2f58d52d66.zh-cn.xlf"
        Handback    = "2016-11-14 06:33:25"
        HyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/37c89f074a5640d30d6b8520b98cd86d7eff4b59/e2e/6daa4734-0ac8-48c1-b198-78425699ee4e.md"
    },
    @{
        SheetName   = "de-de"
        TargetFile  = "6daa4734-0ac8-48c1-b198-78425699ee4e.8db02af582f18e84bbda947aeb8bdd2f58d52d66.de-de.xlf"
        Handback    = "2016-11-14 06:33:43"
        HyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/37c89f074a5640d30d6b8520b98cd86d7eff4b59/e2e/6daa4734-0ac8-48c1-b198-78425699ee4e.md"
    }
)

$ErrorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b26d2fb39d5f58befda7cdda4a8cdbb266091fa2/e2e/6daa4734-0ac8-48c1-b198-78425699ee4e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37c89f074a5640d30d6b8520b98cd86d7eff4b59/e2e/6daa4734-0ac8-48c1-b198-78425699ee4e.md."

foreach ($row in $rows) {
    $ws = $wb.Worksheets.Item($row.SheetName)

    # Widen the Error Detail column (P) so the long message is readable.
    $ws.Columns.Item(16).ColumnWidth = $ColWidth

    # Latest Target File (I6) - hyperlinked, like the rows above it.
    $ws.Hyperlinks.Add($ws.Range("I6"), $row.HyperlinkUrl, "", "", "6daa4734-0ac8-48c1-b198-78425699ee4e.md")

    # Latest Handback File (J6).
    $ws.Range("J6").Value = $row.TargetFile

    # Latest Handback DateTime (K6).
    $ws.Range("K6").Value = $row.Handback

    # Error Detail (P6).
    $ws.Range("P6").Value = $ErrorDetail
}
